# Apply updated crypto price (col D) / 1h-volume (col E) figures, rows 2-51.
# Price values that look numeric get a leading apostrophe so PowerShell/Excel
# store them as literal text (''X' -> cell text is 'X, i.e. apostrophe + X),
# matching the original text-formatted cells instead of turning them into numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.870.74'
$ws.Range("E2").Value = '  -0.91%  '
$ws.Range("D3").Value = '1.859.92'
$ws.Range("E3").Value = '  -0.44%  '
$ws.Range("D4").Value = '''1.000'
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '''304.85'
$ws.Range("E5").Value = '  -0.75%  '
$ws.Range("E6").Value = '  +0.08%  '
$ws.Range("D7").Value = '''0.5048'
$ws.Range("E7").Value = '  -1.13%  '
$ws.Range("D8").Value = '''0.3640'
$ws.Range("E8").Value = '  -2.63%  '
$ws.Range("D9").Value = '''0.07182'
$ws.Range("E9").Value = '  +0.66%  '
$ws.Range("D10").Value = '''0.8949'
$ws.Range("E10").Value = '  +0.76%  '
$ws.Range("D11").Value = '''20.70'
$ws.Range("E11").Value = '  +0.34%  '
$ws.Range("D12").Value = '1.863.81'
$ws.Range("E12").Value = '  +0.00%  '
$ws.Range("D13").Value = '''0.07486'
$ws.Range("E13").Value = '  -0.63%  '
$ws.Range("D14").Value = '''92.40'
$ws.Range("E14").Value = '  +3.49%  '
$ws.Range("D15").Value = '''5.226'
$ws.Range("E15").Value = '  -1.71%  '
$ws.Range("E16").Value = '  +0.37%  '
$ws.Range("D17").Value = '''0.000008467'
$ws.Range("E17").Value = '  +0.11%  '
$ws.Range("E18").Value = '  -0.03%  '
$ws.Range("E19").Value = '  +0.01%  '
$ws.Range("D20").Value = '26.908.22'
$ws.Range("E20").Value = '  -0.91%  '
$ws.Range("E21").Value = '  -0.58%  '
$ws.Range("D22").Value = '2.094.80'
$ws.Range("E22").Value = '  -0.07%  '
$ws.Range("E23").Value = '  -1.99%  '
$ws.Range("D24").Value = '''6.400'
$ws.Range("E24").Value = '  -1.34%  '
$ws.Range("D25").Value = '''147.87'
$ws.Range("E25").Value = '  -1.19%  '
$ws.Range("E26").Value = '  -2.77%  '
$ws.Range("D27").Value = '''17.88'
$ws.Range("E27").Value = '  -0.27%  '
$ws.Range("D28").Value = '''2.060'
$ws.Range("E28").Value = '  -2.03%  '
$ws.Range("D29").Value = '''113.12'
$ws.Range("E29").Value = '  +0.41%  '
$ws.Range("D30").Value = '''4.681'
$ws.Range("E30").Value = '  -1.47%  '
$ws.Range("D31").Value = '''4.672'
$ws.Range("E31").Value = '  -0.22%  '
$ws.Range("D32").Value = '''0.09255'
$ws.Range("E32").Value = '  +2.42%  '
$ws.Range("D33").Value = '''0.05092'
$ws.Range("E33").Value = '  -0.80%  '
$ws.Range("D34").Value = '''0.7446'
$ws.Range("E34").Value = '  +0.95%  '
$ws.Range("D35").Value = '''2.947'
$ws.Range("E35").Value = '  -4.72%  '
$ws.Range("D36").Value = '''1.149'
$ws.Range("E36").Value = '  -1.07%  '
$ws.Range("D37").Value = '''3.285'
$ws.Range("E37").Value = '  +7.92%  '
$ws.Range("D38").Value = '''0.01999'
$ws.Range("E38").Value = '  -2.22%  '
$ws.Range("E39").Value = '  -0.14%  '
$ws.Range("D40").Value = '''0.5519'
$ws.Range("E40").Value = '  +3.77%  '
$ws.Range("D41").Value = '''1.071'
$ws.Range("E41").Value = '  -0.56%  '
$ws.Range("D42").Value = '''118.49'
$ws.Range("E42").Value = '  +1.76%  '
$ws.Range("D43").Value = '''6.483'
$ws.Range("E43").Value = '  -1.85%  '
$ws.Range("D44").Value = '''8.479'
$ws.Range("E44").Value = '  +1.59%  '
$ws.Range("D45").Value = '''0.1469'
$ws.Range("E45").Value = '  -0.29%  '
$ws.Range("D46").Value = '''0.4678'
$ws.Range("E46").Value = '  +0.93%  '
$ws.Range("D47").Value = '''1.000'
$ws.Range("E47").Value = '  +0.09%  '
$ws.Range("D48").Value = '''10.02'
$ws.Range("E48").Value = '  +0.11%  '
$ws.Range("D49").Value = '''1.564'
$ws.Range("E49").Value = '  -0.35%  '
$ws.Range("D50").Value = '''36.98'
$ws.Range("E50").Value = '  +1.26%  '
$ws.Range("D51").Value = '''63.06'
$ws.Range("E51").Value = '  -2.23%  '
